$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1028
$ws1.Range("F4").Value = 167
$ws1.Range("F5").Value = 2780
$ws1.Range("F6").Value = 93
$ws1.Range("F7").Value = 219
$ws1.Range("F10").Value = 61
$ws1.Range("F11").Value = 71
$ws1.Range("F12").Value = 2591
$ws1.Range("F13").Value = 759

# Sheet "全部类型" (all types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1028
$ws4.Range("F5").Value = 167
$ws4.Range("F6").Value = 2780
$ws4.Range("F7").Value = 93
$ws4.Range("F8").Value = 219
$ws4.Range("F12").Value = 61
$ws4.Range("F13").Value = 71
$ws4.Range("F14").Value = 2591
$ws4.Range("F15").Value = 759
